$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.358.52'
$ws.Range('D2').Style = 'Normal'

$ws.Range('E2').Value = '  -3.32%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.678.70'
$ws.Range('D3').Style = 'Normal'

$ws.Range('E3').Value = '  -3.44%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.31'
$ws.Range('D5').Style = 'Normal'

$ws.Range('E5').Value = '  -2.70%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.58'
$ws.Range('D6').Style = 'Normal'

$ws.Range('E6').Value = '  -6.15%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.676.44'
$ws.Range('D7').Style = 'Normal'

$ws.Range('E7').Value = '  -3.44%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.525'
$ws.Range('D9').Style = 'Normal'

$ws.Range('E9').Value = '  -0.26%  '

$ws.Range('E10').Value = '  -4.73%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.18'
$ws.Range('D11').Style = 'Normal'

$ws.Range('E11').Value = '  -4.59%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.462'
$ws.Range('D12').Style = 'Normal'

$ws.Range('E12').Value = '  -4.07%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.48'
$ws.Range('D13').Style = 'Normal'

$ws.Range('E13').Value = '  -6.87%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000240'
$ws.Range('D14').Style = 'Normal'

$ws.Range('E14').Value = '  -5.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.284.59'
$ws.Range('D15').Style = 'Normal'

$ws.Range('E15').Value = '  -3.63%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.673.61'
$ws.Range('D16').Style = 'Normal'

$ws.Range('E16').Value = '  -3.67%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.390.21'
$ws.Range('D17').Style = 'Normal'

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.64'
$ws.Range('D18').Style = 'Normal'

$ws.Range('E18').Value = '  +6.21%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.16'
$ws.Range('D19').Style = 'Normal'

$ws.Range('E19').Value = '  -4.86%  '

$ws.Range('E20').Value = '  -3.64%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '491.24'
$ws.Range('D21').Style = 'Normal'

$ws.Range('E21').Value = '  -2.97%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.13'
$ws.Range('D22').Style = 'Normal'

$ws.Range('E22').Value = '  -4.20%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.720'
$ws.Range('D23').Style = 'Normal'

$ws.Range('E23').Value = '  -2.42%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '85.34'
$ws.Range('D24').Style = 'Normal'

$ws.Range('E24').Value = '  -0.88%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.30'
$ws.Range('D25').Style = 'Normal'

$ws.Range('E25').Value = '  -6.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000137'
$ws.Range('D26').Style = 'Normal'

$ws.Range('E26').Value = '  -3.93%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.14'
$ws.Range('D27').Style = 'Normal'

$ws.Range('E27').Value = '  -3.91%  '

$ws.Range('E28').Value = '  +0.13%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.99'
$ws.Range('D29').Style = 'Normal'

$ws.Range('E29').Value = '  -4.33%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.36'
$ws.Range('D31').Style = 'Normal'

$ws.Range('E31').Value = '  -6.35%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.65'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E32').Value = '  -3.81%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.44'
$ws.Range('D33').Style = 'Normal'

$ws.Range('E33').Value = '  -0.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.811.82'
$ws.Range('D34').Style = 'Normal'

$ws.Range('E34').Value = '  -3.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.108'
$ws.Range('D35').Style = 'Normal'

$ws.Range('E35').Value = '  -4.88%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.611.88'
$ws.Range('D36').Style = 'Normal'

$ws.Range('E36').Value = '  -3.54%  '

$ws.Range('E37').Value = '  -0.03%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.993'
$ws.Range('D38').Style = 'Normal'

$ws.Range('E38').Value = '  -4.59%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.75'
$ws.Range('D39').Style = 'Normal'

$ws.Range('E39').Value = '  -5.53%  '

$ws.Range('E40').Value = '  -7.00%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.323'
$ws.Range('D41').Style = 'Normal'

$ws.Range('E41').Value = '  -3.69%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '432.45'
$ws.Range('D42').Style = 'Normal'

$ws.Range('E42').Value = '  -10.64%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.54'
$ws.Range('D43').Style = 'Normal'

$ws.Range('E43').Value = '  -2.35%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.93'
$ws.Range('D44').Style = 'Normal'

$ws.Range('E44').Value = '  -5.74%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.77'
$ws.Range('D45').Style = 'Normal'

$ws.Range('E45').Value = '  -8.02%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.36'
$ws.Range('D46').Style = 'Normal'

$ws.Range('E46').Value = '  -1.88%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '40.56'
$ws.Range('D48').Style = 'Normal'

$ws.Range('E48').Value = '  -6.97%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.11'
$ws.Range('D49').Style = 'Normal'

$ws.Range('E49').Value = '  +1.83%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.748.95'
$ws.Range('D50').Style = 'Normal'

$ws.Range('E50').Value = '  -5.76%  '

$ws.Range('E51').Value = '  -3.81%  '
